# hash & equals von Order & User gekürzt
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 6 corresponds to order_id = 5 (firm "Weiler")
# Columns: A=order_id B=firm C=stone_type D=amount E=due_date
#          F=phase G=price H=done I=status J=rowcount
$ws.Range("C6").Value = "Sandstein"
$ws.Range("F6").Value = "Planung"
$ws.Range("G6").Value = 60000
